# fix(2025/tempate): update color of the accent bar
#
# The slide master contains a thin vertical "accent bar" rectangle shape
# (named "矩形 2") positioned just to the left of the title placeholder.
# Its solid fill color is updated from the old dark red (CA463A) to the
# new brighter red (FF3535) used elsewhere in the template (e.g. the
# quote-mark text boxes already use FF3535).

$p = $ppt.ActivePresentation
$m = $p.SlideMaster

$bar = $null
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $shp = $m.Shapes.Item($i)
    if ($shp.Name -eq "矩形 2") {
        $bar = $shp
        break
    }
}

if ($bar -eq $null) {
    throw "accent bar shape '矩形 2' not found on slide master"
}

# FF3535 as a PowerPoint BGR-packed RGB long: (0x35 << 16) | (0x35 << 8) | 0xFF
$bar.Fill.ForeColor.RGB = 3487231
